$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '50.938.22'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '2.945.80'
$ws.Range('E3').Value = '  -0.55%  '
$r = $ws.Range('D4')
$r.NumberFormat = '@'
$r.Value = '0.999'
$ws.Range('E4').Value = '  +0.01%  '
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '379.01'
$ws.Range('E5').Value = '  -1.00%  '
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '101.07'
$ws.Range('E6').Value = '  -2.06%  '
$r = $ws.Range('D7')
$r.NumberFormat = '@'
$r.Value = '0.541'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -1.48%  '
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '36.10'
$ws.Range('E10').Value = '  -1.52%  '
$ws.Range('E11').Value = '  -0.61%  '
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '0.0847'
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').Value = '3.403.78'
$ws.Range('E13').Value = '  -0.72%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '7.76'
$ws.Range('E14').Value = '  +3.99%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '18.27'
$ws.Range('E15').Value = '  +1.32%  '
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '12.07'
$ws.Range('E16').Value = '  +68.85%  '
$ws.Range('D17').Value = '2.945.42'
$ws.Range('E17').Value = '  -0.36%  '
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '1.00'
$ws.Range('E18').Value = '  +0.91%  '
$ws.Range('D19').Value = '50.927.82'
$ws.Range('E20').Value = '  -4.56%  '
$ws.Range('E21').Value = '  -1.86%  '
$ws.Range('E22').Value = '  -0.65%  '
$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '69.41'
$ws.Range('E23').Value = '  +1.29%  '
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '266.21'
$ws.Range('E24').Value = '  +1.43%  '
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '3.19'
$ws.Range('E25').Value = '  +9.06%  '
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '8.12'
$ws.Range('E26').Value = '  -3.23%  '
$ws.Range('E27').Value = '  +0.01%  '
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '7.05'
$ws.Range('E28').Value = '  -9.86%  '
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '25.56'
$ws.Range('E29').Value = '  -0.68%  '
$ws.Range('E30').Value = '  -4.15%  '
$ws.Range('E31').Value = '  -4.02%  '
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '10.10'
$ws.Range('E32').Value = '  +2.69%  '
$ws.Range('E33').Value = '  -0.07%  '
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '50.43'
$ws.Range('E34').Value = '  -0.07%  '
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '33.41'
$ws.Range('E35').Value = '  -1.59%  '
$ws.Range('E36').Value = '  -6.17%  '
$ws.Range('E37').Value = '  -0.08%  '
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '3.09'
$ws.Range('E38').Value = '  +3.09%  '
$ws.Range('E39').Value = '  +0.25%  '
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '16.57'
$ws.Range('E40').Value = '  -1.67%  '
$ws.Range('E41').Value = '  +0.95%  '
$ws.Range('E42').Value = '  -1.34%  '
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '118.89'
$ws.Range('E43').Value = '  -2.05%  '
$ws.Range('E44').Value = '  +8.22%  '
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '21.32'
$ws.Range('E45').Value = '  -0.71%  '
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '2.01'
$ws.Range('E46').Value = '  -1.80%  '
$ws.Range('D48').Value = '2.000.60'
$ws.Range('E48').Value = '  -0.71%  '
$ws.Range('E49').Value = '  -5.02%  '
$ws.Range('E50').Value = '  -10.26%  '
$ws.Range('E51').Value = '  +3.78%  '
